$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 138
$ws.Range("F4").Value = 9266
$ws.Range("F7").Value = 6395
$ws.Range("F10").Value = 9768
$ws.Range("F11").Value = 11099
$ws.Range("F13").Value = 1145
$ws.Range("F14").Value = 4919
$ws.Range("F15").Value = 791
$ws.Range("F16").Value = 451
$ws.Range("F18").Value = 331
$ws.Range("F21").Value = 240
$ws.Range("F23").Value = 879
$ws.Range("F24").Value = 1236
$ws.Range("F25").Value = 855
$ws.Range("F27").Value = 2022
$ws.Range("F30").Value = 2653
$ws.Range("F32").Value = 1739
$ws.Range("F33").Value = 93
$ws.Range("F34").Value = 797
$ws.Range("F35").Value = 47
$ws.Range("F37").Value = 590
$ws.Range("F38").Value = 20
$ws.Range("F39").Value = 3311
$ws.Range("F42").Value = 509
$ws.Range("F48").Value = 4204
$ws.Range("F49").Value = 27

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 30
$ws.Range("F23").Value = 68

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5902

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 138
$ws.Range("F4").Value = 9266
$ws.Range("F9").Value = 6395
$ws.Range("F11").Value = 9768
$ws.Range("F12").Value = 11099
$ws.Range("F14").Value = 1145
$ws.Range("F15").Value = 4919
$ws.Range("F16").Value = 791
$ws.Range("F17").Value = 451
$ws.Range("F19").Value = 331
$ws.Range("F23").Value = 240
$ws.Range("F24").Value = 855
$ws.Range("F26").Value = 2022
$ws.Range("F29").Value = 2653
$ws.Range("F31").Value = 1739
$ws.Range("F32").Value = 93
$ws.Range("F34").Value = 797
$ws.Range("F39").Value = 47
$ws.Range("F41").Value = 590
$ws.Range("F42").Value = 20
$ws.Range("F49").Value = 4204
